$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New plain-text cell C1
$ws.Range("C1").Value = "C1 < C2 > C1"

# Rich-text cell: "no format " + bold "bold format" + " no format"
$ws.Range("C2").Value = "no format bold format no format"
$boldRun = $ws.Range("C2").Characters(11, 11).Font
$boldRun.Bold = $true
$tailRun = $ws.Range("C2").Characters(22, 10).Font
$tailRun.Name = "Calibri"
$tailRun.Size = 11

# New plain-text cell D1
$ws.Range("D1").Value = "D1"

# Match the recorded selection / active sheet state
$ws.Range("E9").Select()
